$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# The "Arbitrary Precision" solution row (row 22 on the Data sheet) is
# removed from the report. It was an almost-empty row (no measurements
# were ever filled in), so the row is deleted outright and everything
# below it (the "AXI" row) shifts up to take its place.
# -----------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
$data.Rows("22:22").Delete()

# -----------------------------------------------------------------------
# Every chart that plots the "Data" table references the solution-name
# column (B5:B23) and a measurement column (5:23) as its category/value
# ranges. Re-point each series formula at the new, one-row-shorter range
# (5:22) so the charts stay in sync with the shifted data.
# -----------------------------------------------------------------------
function Set-SeriesFormulas($chart, $cols) {
    for ($i = 0; $i -lt $chart.SeriesCollection().Count; $i++) {
        $col = $cols[$i]
        $chart.SeriesCollection().Item($i + 1).Formula = `
            "=SERIES(Data!`$$col`$3,Data!`$B`$5:`$B`$22,Data!`$$col`$5:`$$col`$22,1)"
    }
}

$dynPlots = $wb.Worksheets.Item("Dynamic Power Plots")
Set-SeriesFormulas $dynPlots.ChartObjects().Item(1).Chart @("D", "E", "F", "G", "H", "I", "J", "K")
Set-SeriesFormulas $dynPlots.ChartObjects().Item(2).Chart @("M")

$utilPlots = $wb.Worksheets.Item("Utilization Report Plots")
Set-SeriesFormulas $utilPlots.ChartObjects().Item(1).Chart @("O", "P", "Q", "R", "S", "T", "U")

$timePlots = $wb.Worksheets.Item("Timing Report Plots")
Set-SeriesFormulas $timePlots.ChartObjects().Item(1).Chart @("W")
Set-SeriesFormulas $timePlots.ChartObjects().Item(2).Chart @("X")
Set-SeriesFormulas $timePlots.ChartObjects().Item(3).Chart @("Y")
Set-SeriesFormulas $timePlots.ChartObjects().Item(4).Chart @("Z")

# -----------------------------------------------------------------------
# Restore the cursor to where the author left it before saving.
# -----------------------------------------------------------------------
$data.Range("B30").Select() | Out-Null
